$d = $word.ActiveDocument

$replacements = @(
    @("304÷7=", "501÷9="),
    @("113÷6=", "597÷2="),
    @("263÷8=", "745÷8="),
    @("815÷4=", "340÷8="),
    @("311÷3=", "318÷5="),
    @("481÷3=", "814÷6="),
    @("943÷4=", "848÷7="),
    @("585÷7=", "113÷8="),
    @("463÷4=", "995÷3="),
    @("906÷6=", "583÷4="),
    @("495÷5=", "880÷3="),
    @("176÷9=", "288÷4="),
    @("909÷5=", "275÷5="),
    @("125÷4=", "397÷5="),
    @("980÷2=", "620÷4="),
    @("668÷5=", "802÷7="),
    @("821÷7=", "582÷5="),
    @("455÷8=", "431÷9="),
    @("422÷7=", "782÷6="),
    @("892÷4=", "761÷8="),
    @("224÷7=", "667÷5="),
    @("716÷4=", "229÷6="),
    @("704÷6=", "910÷5="),
    @("568÷5=", "666÷4="),
    @("320÷8=", "152÷4=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
